$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'34.815.54"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = "'1.866.68"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('E4').Value = '  -0.97%  '
$ws.Range('D5').Value = "'243.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.21%  '
$ws.Range('D6').Value = "'0.669"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.81%  '
$ws.Range('E7').Value = '  -1.00%  '
$ws.Range('D8').Value = "'42.07"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.08%  '
$ws.Range('E9').Value = '  -5.61%  '
$ws.Range('D10').Value = "'0.0732"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').Value = "'0.0967"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').Value = "'12.78"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').Value = "'2.136.85"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.42%  '
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').Value = "'1.866.11"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('E16').Value = '  -2.50%  '
$ws.Range('D17').Value = "'34.781.35"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = "'71.92"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.24%  '
$ws.Range('D19').Value = "'0.0₃0807"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.96%  '
$ws.Range('D20').Value = "'241.87"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('E21').Value = '  -3.81%  '
$ws.Range('E22').Value = '  -4.17%  '
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('E24').Value = '  +5.15%  '
$ws.Range('D25').Value = "'2.14"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -13.23%  '
$ws.Range('D26').Value = "'162.84"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('D27').Value = "'8.28"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.16%  '
$ws.Range('D28').Value = "'17.98"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.94%  '
$ws.Range('E29').Value = '  -5.86%  '
$ws.Range('D30').Value = "'4.128.39"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('D32').Value = "'4.13"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.04%  '
$ws.Range('D33').Value = "'0.0567"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.81%  '
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('D35').Value = "'4.10"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('D36').Value = "'0.826"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.41%  '
$ws.Range('D37').Value = "'1.92"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('D38').Value = "'1.49"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -25.45%  '
$ws.Range('D39').Value = "'97.01"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').Value = "'16.87"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.39%  '
$ws.Range('D41').Value = "'0.0664"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.40%  '
$ws.Range('D42').Value = "'0.0209"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.81%  '
$ws.Range('E43').Value = '  -4.40%  '
$ws.Range('D44').Value = "'0.0826"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +11.81%  '
$ws.Range('D45').Value = "'1.277.24"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.75%  '
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('E47').Value = '  -1.26%  '
$ws.Range('E48').Value = '  -1.71%  '
$ws.Range('D49').Value = "'11.87"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.04%  '
$ws.Range('E50').Value = '  -7.77%  '
$ws.Range('D51').Value = "'42.27"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.49%  '
